# Generate Report for Handoff
# Adds a new tracked file (d08ad933-095b-40df-b7c4-daa762fa9ddc.md) as row 9
# to each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns: File Name | Path And Name |
# Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bba7348ffe03113c13ca99620c42166a1a526839/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md",
    "",
    "",
    "e2e\d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
) | Out-Null
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-18 12:44:37"
$wsOverview.Range("G9").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - columns: Source File Name | File Extension |
# Status | Source Path | Priority | Content Duplicate | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback
# File | Latest Handback DateTime | Reference Tokens | To be localized |
# Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bba7348ffe03113c13ca99620c42166a1a526839/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md",
    "",
    "",
    "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
) | Out-Null
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-18 12:44:31"
$wsZhCn.Range("H9").NumberFormat = $dateFmt
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = $dateFmt
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) - same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bba7348ffe03113c13ca99620c42166a1a526839/e2e/d08ad933-095b-40df-b7c4-daa762fa9ddc.md",
    "",
    "",
    "d08ad933-095b-40df-b7c4-daa762fa9ddc.md"
) | Out-Null
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "d08ad933-095b-40df-b7c4-daa762fa9ddc.bba7348ffe03113c13ca99620c42166a1a526839.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-18 12:44:37"
$wsDeDe.Range("H9").NumberFormat = $dateFmt
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = $dateFmt
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

Write-Host "Added row 9 (d08ad933-095b-40df-b7c4-daa762fa9ddc.md) to Overview, zh-cn, de-de sheets."
